$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (correlation) values
$ws.Range("B2").Value = -0.3400423899126752
$ws.Range("C2").Value = 0.2539039686601514
$ws.Range("D2").Value = -0.1435109388079117
$ws.Range("E2").Value = -0.06402795731429904
$ws.Range("G2").Value = -0.05453510026901395
$ws.Range("H2").Value = 0.2063054549602106
$ws.Range("I2").Value = -0.3164322054900608
$ws.Range("J2").Value = 0.4534878842403504
$ws.Range("K2").Value = 0.08019867686619699
$ws.Range("M2").Value = 0.08759447578528413
$ws.Range("N2").Value = 0.3079629191661964
$ws.Range("O2").Value = 0.3175867603901401
$ws.Range("P2").Value = -0.2671511322424202
$ws.Range("Q2").Value = -0.08226004208754302
$ws.Range("S2").Value = -0.5135060154474426
$ws.Range("T2").Value = 0.3653241690774013
$ws.Range("U2").Value = 0.012117414870467
$ws.Range("V2").Value = -0.1959289114047993

# Row 3 (p-value) values
$ws.Range("B3").Value = 0.007854267608780977
$ws.Range("C3").Value = 0.05028055808055044
$ws.Range("D3").Value = 0.2739896320923436
$ws.Range("E3").Value = 0.6269496924153324
$ws.Range("G3").Value = 0.6789846176908165
$ws.Range("H3").Value = 0.1137673362698121
$ws.Range("I3").Value = 0.01377161471086734
$ws.Range("J3").Value = 0.0002738789944284008
$ws.Range("K3").Value = 0.5424377485709881
$ws.Range("M3").Value = 0.505723161097994
$ws.Range("N3").Value = 0.01667550495892327
$ws.Range("O3").Value = 0.01341161308485215
$ws.Range("P3").Value = 0.03906694655211305
$ws.Range("Q3").Value = 0.5320760097623757
$ws.Range("S3").Value = 0.00002718118566890267
$ws.Range("T3").Value = 0.0000004592859788038711
$ws.Range("U3").Value = 0.8717425259125841
$ws.Range("V3").Value = 0.008390272948327417
